# Weekly Fruit/Vegetable price update:
# Insert a new observation row before the current row 268 (pushing all
# subsequent rows down by one), matching the pattern already used for the
# rest of the sheet (columns A,B,C,E,F,G,R stay constant for this market/
# category combination).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(268).Insert()

$ws.Range("A268").Value = 11
$ws.Range("B268").Value = "Vega Monumental Concepción"
$ws.Range("C268").Value = "Bíobío"
$ws.Range("D268").Value = 44516
$ws.Range("E268").Value = 8
$ws.Range("F268").Value = 100112020
$ws.Range("G268").Value = "Tomate"
$ws.Range("H268").Value = "Larga vida"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 550
$ws.Range("K268").Value = 7500
$ws.Range("L268").Value = 8000
$ws.Range("M268").Value = 7727
$ws.Range("N268").Value = "`$/caja 15 kilos"
$ws.Range("O268").Value = "Región del Maule"
$ws.Range("P268").Value = 515
$ws.Range("Q268").Value = 15
$ws.Range("R268").Value = "Hortaliza"
